# Appends newly scraped 7-ELEVEN store rows (districts belonging to the
# remainder of the county) to the bottom of the sheet, as if "next page"
# of the scraper had just been consumed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = "和豐"
$ws.Range("D7").Value = "基隆市中正區新豐街203號1樓"
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("Q7").Value = 1
$ws.Range("AB7").Value = 1
$ws.Range("AC7").Value = 1
$ws.Range("AI7").Value = 1
$ws.Range("AK7").Value = 1

$ws.Range("C8").Value = "哨船頭"
$ws.Range("D8").Value = "基隆市中正區義一路43號1樓"
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 1
$ws.Range("O8").Value = 1
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = 1
$ws.Range("T8").Value = 1
$ws.Range("AB8").Value = 1
$ws.Range("AC8").Value = 1
$ws.Range("AD8").Value = 1
$ws.Range("AI8").Value = 1

$ws.Range("C9").Value = "海洋"
$ws.Range("D9").Value = "基隆市中正區中正路609號1樓"
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = 1
$ws.Range("S9").Value = 1
$ws.Range("T9").Value = 1
$ws.Range("U9").Value = 1
$ws.Range("Z9").Value = 1
$ws.Range("AB9").Value = 1
$ws.Range("AC9").Value = 1
$ws.Range("AD9").Value = 1
$ws.Range("AF9").Value = 1
$ws.Range("AH9").Value = 1
$ws.Range("AI9").Value = 1
$ws.Range("AJ9").Value = 1
$ws.Range("AK9").Value = 1
$ws.Range("AL9").Value = 1

$ws.Range("C10").Value = "財豐"
$ws.Range("D10").Value = "基隆市中正區新豐街389號1樓"
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = 1
$ws.Range("S10").Value = 1
$ws.Range("AC10").Value = 1
$ws.Range("AD10").Value = 1
$ws.Range("AI10").Value = 1

$ws.Range("C11").Value = "基義"
$ws.Range("D11").Value = "基隆市中正區義二路181號185號1樓"
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("Q11").Value = 1
$ws.Range("S11").Value = 1
$ws.Range("AB11").Value = 1
$ws.Range("AC11").Value = 1
$ws.Range("AD11").Value = 1
$ws.Range("AG11").Value = 1
$ws.Range("AI11").Value = 1

$ws.Range("C12").Value = "港都"
$ws.Range("D12").Value = "基隆市中正區義二路8號1樓"
$ws.Range("G12").Value = 1
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = 1
$ws.Range("S12").Value = 1
$ws.Range("AB12").Value = 1
$ws.Range("AC12").Value = 1
$ws.Range("AD12").Value = 1
$ws.Range("AI12").Value = 1

$ws.Range("C13").Value = "翔濱"
$ws.Range("D13").Value = "基隆市中正區建國里祥豐街339號"
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 1
$ws.Range("L13").Value = 1
$ws.Range("Q13").Value = 1
$ws.Range("AB13").Value = 1
$ws.Range("AC13").Value = 1

$ws.Range("C14").Value = "新財發"
$ws.Range("D14").Value = "基隆市中正區新豐街303巷11弄1號3號1樓"
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 1
$ws.Range("N14").Value = 1
$ws.Range("Q14").Value = 1
$ws.Range("AB14").Value = 1
$ws.Range("AC14").Value = 1
$ws.Range("AI14").Value = 1
$ws.Range("AK14").Value = 1

$ws.Range("C15").Value = "漁港"
$ws.Range("D15").Value = "基隆市中正區中正路672號"
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 1
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 1
$ws.Range("Q15").Value = 1
$ws.Range("AB15").Value = 1
$ws.Range("AC15").Value = 1
$ws.Range("AD15").Value = 1
$ws.Range("AF15").Value = 1
$ws.Range("AI15").Value = 1

$ws.Range("C16").Value = "龍騰"
$ws.Range("D16").Value = "基隆市中正區義一路22號24號"
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 1
$ws.Range("L16").Value = 1
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = 1
$ws.Range("S16").Value = 1
$ws.Range("T16").Value = 1
$ws.Range("AB16").Value = 1
$ws.Range("AC16").Value = 1
$ws.Range("AI16").Value = 1

$ws.Range("C17").Value = "豐勝"
$ws.Range("D17").Value = "基隆市中正區中正路322號"
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = 1
$ws.Range("Q17").Value = 1
$ws.Range("AB17").Value = 1
$ws.Range("AC17").Value = 1
$ws.Range("AI17").Value = 1

$ws.Range("C18").Value = "觀山海"
$ws.Range("D18").Value = "基隆市中正區砂子里觀海街49號51號1樓"
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("Q18").Value = 1
$ws.Range("S18").Value = 1
$ws.Range("AB18").Value = 1
$ws.Range("AC18").Value = 1
$ws.Range("AD18").Value = 1
$ws.Range("AI18").Value = 1
$ws.Range("AK18").Value = 1

